$ErrorActionPreference = "Stop"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Protect()
$ws.Unprotect()
$ws.Range("F1").Value = "Pengganti"
Write-Output "set ok"
